$d = $word.ActiveDocument

# Collapse to the very end of the document body (right after the
# existing "Local Nov 7 file" paragraph).
$endRange = $d.Content
$endRange.Collapse(0)   # wdCollapseEnd

# Build the two new paragraphs as raw OOXML so we get exactly:
#   1) a completely empty paragraph (<w:p/>)
#   2) a paragraph containing two separate runs of text
#      ("Dev 2 starts working. Congrats" and ". You did it.")
# Doing this as a single InsertXML call keeps both paragraphs distinct
# (two separate InsertXML/InsertAfter calls at "end of document" would
# instead land inside / merge with the previously inserted paragraph,
# since that position is simultaneously "end of that paragraph" and
# "end of the document").
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParagraphsXml = "<w:p $wNs/>" +
                     "<w:p $wNs>" +
                       "<w:r><w:t>Dev 2 starts working. Congrats</w:t></w:r>" +
                       "<w:r><w:t>. You did it.</w:t></w:r>" +
                     "</w:p>"

[void]$endRange.InsertXML($newParagraphsXml)
